$wb = $excel.ActiveWorkbook

# Rename the "total_staff" sheet to "total_staff_by_type"
$ws = $wb.Worksheets.Item("total_staff")
$ws.Name = "total_staff_by_type"

# Update the active cell selection on that sheet from D2 to D17
$ws.Activate()
$ws.Range("D17").Select()
